$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after "City of Frederick"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Town of Lonaconing"

# Header row
$ws2.Range("A1").Value = "Zone"
$ws2.Range("B1").Value = "Zone Abbreviation"
$ws2.Range("C1").Value = "Issuing Body"
$ws2.Range("D1").Value = "Zone General Description"

# Row 2 - Low Density Residential (first pass: abbreviation + issuing body)
$ws2.Range("A2").Value = "Low Density Residential"
$ws2.Range("B2").Value = "LDR"
$ws2.Range("C2").Value = "Town of Lonaconing Planning and Zoning Department"

# Row 4 - Town Center
$ws2.Range("B4").Value = "TC"
$ws2.Range("D4").Value = "The Town Center (TC) Zoning District is comprised of single family detached dwellings, one dwelling unit in combination with permitted commercial use, twin dwelling, two family detached dwelling, townhouse, and multifamily dwelling."
$ws2.Range("A4").Value = "Town Center"
$ws2.Range("C4").Value = "Town of Lonaconing Planning and Zoning Department"

# Row 3 - Medium Density Residential
$ws2.Range("A3").Value = "Medium Density Residential"
$ws2.Range("B3").Value = "MDR"
$ws2.Range("D3").Value = "The Medium Density Residential (MDR) Zoning District is comprised of single family detached dwilling, twin dwelling, two family detached dwelling, townhouse, and multifamily dwelling."
$ws2.Range("C3").Value = "Town of Lonaconing Planning and Zoning Department"

# Row 5 - Mixed Use
$ws2.Range("A5").Value = "Mixed Use "
$ws2.Range("B5").Value = "MU"
$ws2.Range("D5").Value = "The Mixed Use (MU) Zoning District is any structure or use; however any dwelling units should include a minimum of 500 square feet of indoor floor area per dwelling unit."
$ws2.Range("C5").Value = "Town of Lonaconing Planning and Zoning Department"

# Row 2 - Low Density Residential (final piece: general description, added last)
$ws2.Range("D2").Value = "The Low Density Residential(LDR) Zoning District is primarily comprised of single family detached dwellings. Although twin dwellings (side-by-side dwellings) are permitted in the LDR as well."

# Column widths matching the "best fit" sizes from the target workbook
# (closest values achievable given the engine's internal 1/6-character rounding)
$ws2.Columns.Item(1).ColumnWidth = 25
$ws2.Columns.Item(2).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(3).ColumnWidth = 48.333333333333336
$ws2.Columns.Item(4).ColumnWidth = 170.33333333333334

# Selections / active cells to match target view state
$ws1.Range("A1:D1").Select()
$ws2.Range("D4").Select()

# Make the new sheet the active one (it becomes the selected tab)
$ws2.Activate()
